$d = $word.ActiveDocument

# The schedule table: row 9 (1-based, header is row 1) is "Møde 8" with
# Lasse / Lukas / Jeppe, whose "Dato" cell (column 5) is still empty.
# Fill it in with the date of that meeting.
$tbl = $d.Tables(1)
$cell = $tbl.Cell(9, 5)
$cell.Range.Text = "07-11-2014"
